$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Branch 2004 restricted"
$ws.Range("B12").Value = "p.4"
$ws.Range("C12").Value = "negative feedback from expectations. This particular self-referential feature generates expectations driven oscillations."

# Match formatting of the rest of column C (wrap text), mirroring C11's style.
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C13").Select()
